# Apply the cryptos-list refresh (prices + 1h volume %) described by the commit diff.
# Every D/E data cell in this sheet is stored as TEXT (t="inlineStr" / shared string),
# never as a native number, so each write below pins the cell to the "@" (Text) number
# format before assigning the string, then resets the style back to "Normal" so no stray
# cell-format (s=) attribute leaks into the saved sheet. Only the cell VALUE changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "26.932.11"
Set-TextCell "E2" "  +0.06%  "
Set-TextCell "D3" "1.817.30"
Set-TextCell "E3" "  +0.40%  "
Set-TextCell "D4" "1.003"
Set-TextCell "E4" "  +0.25%  "
Set-TextCell "D5" "309.56"
Set-TextCell "E5" "  -0.29%  "
Set-TextCell "E6" "  +0.17%  "
Set-TextCell "D7" "0.4656"
Set-TextCell "E7" "  +0.62%  "
Set-TextCell "D8" "0.3660"
Set-TextCell "E8" "  -1.44%  "
Set-TextCell "D9" "0.07365"
Set-TextCell "E9" "  -0.19%  "
Set-TextCell "D10" "0.8703"
Set-TextCell "E10" "  -0.71%  "
Set-TextCell "E11" "  -1.29%  "
Set-TextCell "D12" "1.834.03"
Set-TextCell "E12" "  +2.88%  "
Set-TextCell "D13" "5.388"
Set-TextCell "E13" "  +0.50%  "
Set-TextCell "D14" "0.07129"
Set-TextCell "E14" "  +1.35%  "
Set-TextCell "D15" "6.507"
Set-TextCell "E15" "  -0.18%  "
Set-TextCell "D16" "91.34"
Set-TextCell "E16" "  -1.18%  "
Set-TextCell "D17" "1.003"
Set-TextCell "E17" "  +0.41%  "
Set-TextCell "D18" "0.000008683"
Set-TextCell "E18" "  -0.25%  "
Set-TextCell "E19" "  +0.12%  "
Set-TextCell "E20" "  -0.71%  "
Set-TextCell "D21" "26.958.96"
Set-TextCell "E21" "  +0.15%  "
Set-TextCell "D22" "5.290"
Set-TextCell "E22" "  -0.65%  "
Set-TextCell "D23" "10.58"
Set-TextCell "E23" "  -0.57%  "
Set-TextCell "D24" "2.048.54"
Set-TextCell "E24" "  +1.66%  "
Set-TextCell "E25" "  -0.17%  "
Set-TextCell "D26" "151.08"
Set-TextCell "E26" "  -0.16%  "
Set-TextCell "D27" "18.43"
Set-TextCell "E27" "  +0.00%  "
Set-TextCell "D28" "2.139"
Set-TextCell "E28" "  -0.75%  "
Set-TextCell "E29" "  -1.82%  "
Set-TextCell "D30" "116.75"
Set-TextCell "E30" "  +0.55%  "
Set-TextCell "D31" "0.08902"
Set-TextCell "E31" "  -0.02%  "
Set-TextCell "D32" "0.7576"
Set-TextCell "E32" "  +0.17%  "
Set-TextCell "E33" "  +0.51%  "
Set-TextCell "E34" "  +0.56%  "
Set-TextCell "D35" "2.900"
Set-TextCell "E35" "  -0.38%  "
Set-TextCell "E36" "  +0.15%  "
Set-TextCell "D37" "1.096"
Set-TextCell "E37" "  -0.74%  "
Set-TextCell "D38" "0.05286"
Set-TextCell "E38" "  +0.74%  "
Set-TextCell "D39" "0.01945"
Set-TextCell "E39" "  -1.42%  "
Set-TextCell "E40" "  +1.82%  "
Set-TextCell "E41" "  -0.84%  "
Set-TextCell "D42" "7.145"
Set-TextCell "E42" "  -1.06%  "
Set-TextCell "D43" "2.333"
Set-TextCell "E43" "  -3.42%  "
Set-TextCell "D44" "0.1657"
Set-TextCell "E44" "  -0.38%  "
Set-TextCell "D45" "8.421"
Set-TextCell "E45" "  -1.13%  "
Set-TextCell "D46" "0.4849"
Set-TextCell "E46" "  -2.85%  "
Set-TextCell "D47" "10.47"
Set-TextCell "E47" "  +1.41%  "
Set-TextCell "E48" "  +0.17%  "
Set-TextCell "D49" "103.28"
Set-TextCell "E49" "  -0.61%  "
Set-TextCell "E50" "  -1.01%  "
Set-TextCell "E51" "  -0.11%  "
